# Auto-generated script to apply market-data refresh changes to Sheets/Behemoth_Profits.xlsx
# Updates columns H-N (price/profit calc columns) across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1265
$ws.Range("I31").Value = 1265
$ws.Range("K31").Value = 3795
$ws.Range("M31").Value = -3565

$ws.Range("H39").Value = 243.58333
$ws.Range("I39").Value = 145.91667
$ws.Range("K39").Value = 437.75001
$ws.Range("M39").Value = -141.75001

$ws.Range("H52").Value = 2016.5
$ws.Range("I52").Value = 420
$ws.Range("K52").Value = 1260
$ws.Range("M52").Value = -1100

$ws.Range("H107").Value = 360.5
$ws.Range("I107").Value = 360.5
$ws.Range("K107").Value = 360.5
$ws.Range("M107").Value = 1559.5

$ws.Range("H127").Value = 4166.65
$ws.Range("I127").Value = 3215.25
$ws.Range("J127").Value = 5593.75
$ws.Range("K127").Value = 9645.75
$ws.Range("L127").Value = 16781.25
$ws.Range("M127").Value = -4685.75
$ws.Range("N127").Value = -26701.25

$ws.Range("H132").Value = 3910.3076
$ws.Range("I132").Value = 3748.9167
$ws.Range("K132").Value = 11246.7501
$ws.Range("M132").Value = -8716.750100000001

$ws.Range("H137").Value = 4519
$ws.Range("I137").Value = 2699.24
$ws.Range("J137").Value = 10205.75
$ws.Range("K137").Value = 8097.719999999999
$ws.Range("L137").Value = 30617.25
$ws.Range("M137").Value = -5547.719999999999
$ws.Range("N137").Value = -35717.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1298.6
$ws.Range("J2").Value = 1250
$ws.Range("L2").Value = 1250
$ws.Range("N2").Value = -1476

$ws.Range("H74").Value = 5959335
$ws.Range("I74").Value = 8931010
$ws.Range("J74").Value = 15984.643
$ws.Range("K74").Value = 8931010
$ws.Range("L74").Value = 15984.643
$ws.Range("M74").Value = -8930136
$ws.Range("N74").Value = -17732.643

$ws.Range("H77").Value = 5959335
$ws.Range("I77").Value = 8931010
$ws.Range("J77").Value = 15984.643
$ws.Range("K77").Value = 44655050
$ws.Range("L77").Value = 79923.215
$ws.Range("M77").Value = -44650682
$ws.Range("N77").Value = -88659.215

$ws.Range("H88").Value = 2120.6287
$ws.Range("I88").Value = 1295.7
$ws.Range("K88").Value = 1295.7
$ws.Range("M88").Value = -889.7

$ws.Range("H91").Value = 2120.6287
$ws.Range("I91").Value = 1295.7
$ws.Range("K91").Value = 1295.7
$ws.Range("M91").Value = 108.3

$ws.Range("H116").Value = 1298.6
$ws.Range("J116").Value = 1250
$ws.Range("L116").Value = 1250
$ws.Range("N116").Value = -5838

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1298.6
$ws.Range("J3").Value = 1250
$ws.Range("L3").Value = 1250
$ws.Range("N3").Value = -1478

$ws.Range("H18").Value = 1254.5
$ws.Range("I18").Value = 9
$ws.Range("J18").Value = 2500
$ws.Range("K18").Value = 9
$ws.Range("L18").Value = 2500
$ws.Range("M18").Value = 520
$ws.Range("N18").Value = -3558

$ws.Range("H21").Value = 65537
$ws.Range("J21").Value = 65537
$ws.Range("L21").Value = 65537
$ws.Range("N21").Value = -66009

$ws.Range("H36").Value = 2066.8
$ws.Range("I36").Value = 2066.8
$ws.Range("K36").Value = 2066.8
$ws.Range("M36").Value = -1532.8

$ws.Range("H134").Value = 22175608
$ws.Range("I134").Value = 956.1579
$ws.Range("K134").Value = 2868.4737
$ws.Range("M134").Value = -333.4737

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 696300.9399999999
$ws.Range("I31").Value = 13610.0625
$ws.Range("J31").Value = 1303137.2
$ws.Range("K31").Value = 13610.0625
$ws.Range("L31").Value = 1303137.2
$ws.Range("M31").Value = -13315.0625
$ws.Range("N31").Value = -1303727.2

$ws.Range("H34").Value = 696300.9399999999
$ws.Range("I34").Value = 13610.0625
$ws.Range("J34").Value = 1303137.2
$ws.Range("K34").Value = 13610.0625
$ws.Range("L34").Value = 1303137.2
$ws.Range("M34").Value = -13408.0625
$ws.Range("N34").Value = -1303541.2

$ws.Range("H37").Value = 20057
$ws.Range("J37").Value = 20057
$ws.Range("L37").Value = 20057
$ws.Range("N37").Value = -20271

$ws.Range("H64").Value = 49995
$ws.Range("J64").Value = 49995
$ws.Range("L64").Value = 49995
$ws.Range("N64").Value = -50491

$ws.Range("H67").Value = 49995
$ws.Range("J67").Value = 49995
$ws.Range("L67").Value = 49995
$ws.Range("N67").Value = -51711

$ws.Range("H88").Value = 17026.666
$ws.Range("J88").Value = 24040
$ws.Range("L88").Value = 24040
$ws.Range("N88").Value = -24852

$ws.Range("H91").Value = 17026.666
$ws.Range("J91").Value = 24040
$ws.Range("L91").Value = 24040
$ws.Range("N91").Value = -26848

$ws.Range("H107").Value = 1764.6111
$ws.Range("J107").Value = 1803.8182
$ws.Range("L107").Value = 1803.8182
$ws.Range("N107").Value = -5643.8182

$ws.Range("H132").Value = 2225.7446
$ws.Range("I132").Value = 2181.1707
$ws.Range("K132").Value = 6543.5121
$ws.Range("M132").Value = -4013.5121

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 612.2
$ws.Range("I103").Value = 221
$ws.Range("K103").Value = 663
$ws.Range("M103").Value = 216

$ws.Range("H128").Value = 452469.75
$ws.Range("I128").Value = 452469.75
$ws.Range("K128").Value = 1357409.25
$ws.Range("M128").Value = -1352429.25

$ws.Range("H131").Value = 4369.906
$ws.Range("J131").Value = 3800.1064
$ws.Range("L131").Value = 11400.3192
$ws.Range("N131").Value = -21480.3192

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 377500
$ws.Range("I19").Value = 750000
$ws.Range("K19").Value = 750000
$ws.Range("M19").Value = -749712

$ws.Range("H59").Value = 7000
$ws.Range("J59").Value = 8000
$ws.Range("L59").Value = 8000
$ws.Range("N59").Value = -9166

$ws.Range("H80").Value = 1600
$ws.Range("J80").Value = 1833.3334
$ws.Range("L80").Value = 1833.3334
$ws.Range("N80").Value = -3829.3334

$ws.Range("H83").Value = 1600
$ws.Range("J83").Value = 1833.3334
$ws.Range("L83").Value = 9166.666999999999
$ws.Range("N83").Value = -19150.667

$ws.Range("H132").Value = 66676650

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3766666.8
$ws.Range("J2").Value = 3766666.8
$ws.Range("L2").Value = 3766666.8
$ws.Range("N2").Value = -3766890.8

$ws.Range("H35").Value = 1105.2222
$ws.Range("I35").Value = 1394.2858
$ws.Range("J35").Value = 93.5
$ws.Range("K35").Value = 1394.2858
$ws.Range("L35").Value = 93.5
$ws.Range("M35").Value = -1058.2858
$ws.Range("N35").Value = -765.5

$ws.Range("H53").Value = 6732.6665
$ws.Range("I53").Value = 6732.6665
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 6732.6665
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -6214.6665
$ws.Range("N53").ClearContents()

$ws.Range("H58").Value = 28414.25
$ws.Range("J58").Value = 46329
$ws.Range("L58").Value = 46329
$ws.Range("N58").Value = -46849

$ws.Range("H115").Value = 110000
$ws.Range("J115").Value = 110000
$ws.Range("L115").Value = 110000
$ws.Range("N115").Value = -112350

$ws.Range("H132").Value = 559015.8
$ws.Range("I132").Value = 3252.3076
$ws.Range("J132").Value = 2004001
$ws.Range("K132").Value = 9756.9228
$ws.Range("L132").Value = 6012003
$ws.Range("M132").Value = -7226.9228
$ws.Range("N132").Value = -6017063

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1728.9
$ws.Range("I136").Value = 1365.4445
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 4096.333500000001
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -1546.333500000001
$ws.Range("N136").Value = -20100
